# Paralel to Seri Revizyon Board için URL uzantısı eklendi.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

$url = "https://github.com/btk42/EQ-RVB-00-PTS-REV-S3B0-01"
$cell = $ws.Range("K7")

# Add the hyperlink for the "Revizyon Board" (Parallel To Serial) row - text and
# target both the repository URL, same pattern already used by K2/K3.
$ws.Hyperlinks.Add($cell, $url)

# Re-apply the existing "hyperlink" cell style (font/border) used by K2/K3 so the
# new cell matches the sheet's established hyperlink formatting exactly.
$ws.Range("K2").Copy()
$cell.PasteSpecial(-4122)

# Leave the selection where the author ended up after making the edit.
$ws.Range("N11").Select() | Out-Null
